{"js": "// The document's text was re-saved through a pipeline that mangled its\n// UTF-8 encoding: characters that were originally encoded correctly as\n// UTF-8 got their UTF-8 bytes re-interpreted as Windows-1252/Latin-1,\n// producing \"mojibake\". Two characters are affected in this document:\n//   EN DASH        U+2013 \"\u2013\"  ->  \"\u00e2\u20ac\u201c\"  (U+00E2 U+20AC U+201C)\n//   BULLET         U+2022 \"\u2022\"  ->  \"\u00e2\u20ac\u00a2\"  (U+00E2 U+20AC U+00A2)\n// Separately, the commit removes the decorative horizontal-rule\n// paragraph (a <w:pict> containing a zero-width v:rect with o:hr=\"t\"),\n// and the paragraph that used to follow it (\"MANUSCRIPT HIGHLIGHTS\u2026\")\n// reverts from the \"First Paragraph\" style (used right after a\n// horizontal rule for spacing) back to plain \"Body Text\".\n\nconst EN_DASH = \"\\u2013\";\nconst BULLET = \"\\u2022\";\nconst MOJIBAKE_EN_DASH = \"\\u00e2\\u20ac\\u201c\";\nconst MOJIBAKE_BULLET = \"\\u00e2\\u20ac\\u00a2\";\n\nconst body = context.document.body;\n\n// --- 1) Replace every EN DASH occurrence with its mojibake form ------\nlet dashResults = body.search(EN_DASH, { matchCase: true });\ndashResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < dashResults.items.length; i++) {\n  dashResults.items[i].insertText(MOJIBAKE_EN_DASH, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Replace every BULLET occurrence with its mojibake form -------\nlet bulletResults = body.search(BULLET, { matchCase: true });\nbulletResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < bulletResults.items.length; i++) {\n  bulletResults.items[i].insertText(MOJIBAKE_BULLET, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3) Remove the horizontal-rule paragraph and fix up the style ----\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text,style\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  // The horizontal-rule paragraph has no text runs at all (its only\n  // run holds a <w:pict> drawing, not a <w:t>), so it is the unique\n  // empty paragraph in this letter.\n  if (para.text === \"\") {\n    const nextPara = i + 1 < paragraphs.items.length ? paragraphs.items[i + 1] : null;\n    if (nextPara) {\n      nextPara.load(\"style\");\n      await context.sync();\n      if (nextPara.style === \"First Paragraph\") {\n        nextPara.style = \"Body Text\";\n      }\n    }\n    para.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document's text was re-saved through a pipeline that mangled its\n# UTF-8 encoding: characters that were originally correct UTF-8 got\n# their bytes re-interpreted as Windows-1252/Latin-1 (\"mojibake\").\n# Two characters are affected throughout this letter:\n#   EN DASH  U+2013 \"-\"  ->  \"a^co\" (U+00E2 U+20AC U+201C) i.e. \"\u00e2\u20ac\u201c\"\n#   BULLET   U+2022 \"*\"  ->  \"a^c.\"  (U+00E2 U+20AC U+00A2) i.e. \"\u00e2\u20ac\u00a2\"\n# (ASCII-art above is just for the comment; the real characters are used\n# in the code below.) Separately, the decorative horizontal-rule\n# paragraph (an empty paragraph holding only a <w:pict> rule) is\n# deleted, and the paragraph that used to follow it (\"MANUSCRIPT\n# HIGHLIGHTS\u2026\") reverts from the \"First Paragraph\" style back to plain\n# \"Body Text\".\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n$bullet = [char]0x2022\n$mojibakeEnDash = [string]([char]0x00E2) + [string]([char]0x20AC) + [string]([char]0x201C)\n$mojibakeBullet = [string]([char]0x00E2) + [string]([char]0x20AC) + [string]([char]0x00A2)\n\n# --- 1) Replace every EN DASH with its mojibake form, document-wide ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $enDash\n$find.Replacement.Text = $mojibakeEnDash\n$find.Wrap = 1\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# --- 2) Replace every BULLET with its mojibake form, document-wide ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = $bullet\n$find2.Replacement.Text = $mojibakeBullet\n$find2.Wrap = 1\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# --- 3) Remove the horizontal-rule paragraph and fix up the next style\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs($i)\n  $t = $p.Range().Text\n  if ($t.Trim().Length -eq 0) {\n    if ($i -lt $count) {\n      $nextPara = $d.Paragraphs($i + 1)\n      if ($nextPara.Style.NameLocal -eq \"First Paragraph\") {\n        $nextPara.Style = \"Body Text\"\n      }\n    }\n    $p.Range().Delete() | Out-Null\n    break\n  }\n}\n"}
